$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" header timestamp (F1)
$ws.Range("F1").Value = "Last status check on: 25.02.2022 15:45"

# Update row 2 (TankONO): D2 delta price becomes numeric, E2 date becomes numeric date serial
$ws.Range("D2").Value = 0.6
$ws.Range("E2").Value = 44617.64586805556
$ws.Range("E2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
